# Auto-generated PowerShell COM-interop script
# Applies numeric corrections to H:N profit columns across sheets per commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1820.3334
$ws.Range("I18").Value = 1864.4
$ws.Range("K18").Value = 1864.4
$ws.Range("M18").Value = -1580.4
$ws.Range("H28").Value = 514.1
$ws.Range("I28").Value = 567.625
$ws.Range("J28").Value = 300
$ws.Range("K28").Value = 567.625
$ws.Range("L28").Value = 300
$ws.Range("M28").Value = -82.625
$ws.Range("N28").Value = -1270
$ws.Range("H43").Value = 4975
$ws.Range("I43").Value = 1700
$ws.Range("K43").Value = 1700
$ws.Range("M43").Value = -1631
$ws.Range("H86").Value = 4017.5334
$ws.Range("I86").Value = 2399.8333
$ws.Range("K86").Value = 2399.8333
$ws.Range("M86").Value = -1276.8333
$ws.Range("H88").Value = 1666.6
$ws.Range("J88").Value = 2010.4286
$ws.Range("L88").Value = 2010.4286
$ws.Range("N88").Value = -2822.4286
$ws.Range("H89").Value = 4017.5334
$ws.Range("I89").Value = 2399.8333
$ws.Range("K89").Value = 11999.1665
$ws.Range("M89").Value = -6383.166499999999
$ws.Range("H91").Value = 1666.6
$ws.Range("J91").Value = 2010.4286
$ws.Range("L91").Value = 2010.4286
$ws.Range("N91").Value = -4818.4286
$ws.Range("H92").Value = 1282.5385
$ws.Range("I92").Value = 1208.25
$ws.Range("J92").Value = 1401.4
$ws.Range("K92").Value = 1208.25
$ws.Range("L92").Value = 1401.4
$ws.Range("M92").Value = 39.75
$ws.Range("N92").Value = -3897.4
$ws.Range("H132").Value = 1048.3684
$ws.Range("I132").Value = 1001.3333
$ws.Range("K132").Value = 3003.9999
$ws.Range("M132").Value = -473.9998999999998
$ws.Range("H137").Value = 2527.6428
$ws.Range("I137").Value = 1133
$ws.Range("J137").Value = 3573.625
$ws.Range("K137").Value = 3399
$ws.Range("L137").Value = 10720.875
$ws.Range("M137").Value = -849
$ws.Range("N137").Value = -15820.875
$ws.Range("H141").Value = 5297.6665
$ws.Range("I141").Value = 5280
$ws.Range("J141").Value = 5333
$ws.Range("K141").Value = 15840
$ws.Range("L141").Value = 15999
$ws.Range("M141").Value = -10660
$ws.Range("N141").Value = -26359

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15401.099
$ws.Range("I32").Value = 6868.206
$ws.Range("J32").Value = 26146.223
$ws.Range("K32").Value = 6868.206
$ws.Range("L32").Value = 26146.223
$ws.Range("M32").Value = -6581.206
$ws.Range("N32").Value = -26720.223
$ws.Range("H97").Value = 0
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("M97").ClearContents()
$ws.Range("N97").ClearContents()
$ws.Range("H122").Value = 718082.2
$ws.Range("I122").Value = 1114239
$ws.Range("K122").Value = 3342717
$ws.Range("M122").Value = -3340267

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 500
$ws.Range("I94").Value = 500
$ws.Range("K94").Value = 500
$ws.Range("M94").Value = -49
$ws.Range("H134").Value = 3472.0417
$ws.Range("I134").Value = 3492.6086
$ws.Range("J134").Value = 2999
$ws.Range("K134").Value = 10477.8258
$ws.Range("L134").Value = 8997
$ws.Range("M134").Value = -7942.825800000001
$ws.Range("N134").Value = -14067

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 657.8889
$ws.Range("I16").Value = 552.625
$ws.Range("J16").Value = 1500
$ws.Range("K16").Value = 552.625
$ws.Range("L16").Value = 1500
$ws.Range("M16").Value = -265.625
$ws.Range("N16").Value = -2074
$ws.Range("H22").Value = 317.23077
$ws.Range("I22").Value = 366
$ws.Range("K22").Value = 366
$ws.Range("M22").Value = -16
$ws.Range("H62").Value = 59454.715
$ws.Range("I62").Value = 2697.3333
$ws.Range("K62").Value = 2697.3333
$ws.Range("M62").Value = -2073.3333
$ws.Range("H65").Value = 59454.715
$ws.Range("I65").Value = 2697.3333
$ws.Range("K65").Value = 13486.6665
$ws.Range("M65").Value = -10366.6665
$ws.Range("H113").Value = 657.8889
$ws.Range("I113").Value = 552.625
$ws.Range("J113").Value = 1500
$ws.Range("K113").Value = 552.625
$ws.Range("L113").Value = 1500
$ws.Range("M113").Value = 1617.375
$ws.Range("N113").Value = -5840

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 865.36365
$ws.Range("I113").Value = 1315.5
$ws.Range("J113").Value = 765.3333
$ws.Range("K113").Value = 3946.5
$ws.Range("L113").Value = 2295.9999
$ws.Range("M113").Value = -1776.5
$ws.Range("N113").Value = -6635.9999
$ws.Range("H131").Value = 4026.2307
$ws.Range("J131").Value = 4631.1665
$ws.Range("L131").Value = 13893.4995
$ws.Range("N131").Value = -23973.4995
$ws.Range("H132").Value = 3820.5557
$ws.Range("I132").Value = 673.25
$ws.Range("J132").Value = 6338.4
$ws.Range("K132").Value = 6059.25
$ws.Range("L132").Value = 57045.6
$ws.Range("M132").Value = -3529.25
$ws.Range("N132").Value = -62105.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 13300
$ws.Range("J15").Value = 13300
$ws.Range("L15").Value = 13300
$ws.Range("N15").Value = -13876
$ws.Range("H80").Value = 4589.909
$ws.Range("I80").Value = 3471.4285
$ws.Range("J80").Value = 6547.25
$ws.Range("K80").Value = 3471.4285
$ws.Range("L80").Value = 6547.25
$ws.Range("M80").Value = -2473.4285
$ws.Range("N80").Value = -8543.25
$ws.Range("H81").Value = 13300
$ws.Range("J81").Value = 13300
$ws.Range("L81").Value = 13300
$ws.Range("N81").Value = -15296
$ws.Range("H83").Value = 4589.909
$ws.Range("I83").Value = 3471.4285
$ws.Range("J83").Value = 6547.25
$ws.Range("K83").Value = 17357.1425
$ws.Range("L83").Value = 32736.25
$ws.Range("M83").Value = -12365.1425
$ws.Range("N83").Value = -42720.25
$ws.Range("H84").Value = 13300
$ws.Range("J84").Value = 13300
$ws.Range("L84").Value = 39900
$ws.Range("N84").Value = -49884
$ws.Range("H136").Value = 25267
$ws.Range("J136").Value = 25267
$ws.Range("L136").Value = 75801
$ws.Range("N136").Value = -80901

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("H55").Value = 340.85715
$ws.Range("J55").Value = 363.5
$ws.Range("L55").Value = 363.5
$ws.Range("N55").Value = -709.5
$ws.Range("H61").Value = 2574.5
$ws.Range("I61").Value = 2168.4614
$ws.Range("K61").Value = 2168.4614
$ws.Range("M61").Value = -1966.4614
$ws.Range("H113").Value = 2574.5
$ws.Range("I113").Value = 2168.4614
$ws.Range("K113").Value = 2168.4614
$ws.Range("M113").Value = 1.53859999999986

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1859.8
$ws.Range("I81").Value = 1999.8
$ws.Range("J81").Value = 1719.8
$ws.Range("K81").Value = 3999.6
$ws.Range("L81").Value = 3439.6
$ws.Range("M81").Value = -2938.6
$ws.Range("N81").Value = -5561.6
$ws.Range("H84").Value = 1859.8
$ws.Range("I84").Value = 1999.8
$ws.Range("J84").Value = 1719.8
$ws.Range("K84").Value = 19998
$ws.Range("L84").Value = 17198
$ws.Range("M84").Value = -14694
$ws.Range("N84").Value = -27806
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()
$ws.Range("H126").Value = 64818
$ws.Range("I126").Value = 72863.42999999999
$ws.Range("K126").Value = 218590.29
$ws.Range("M126").Value = -216120.29
$ws.Range("H136").Value = 60647.766
$ws.Range("I136").Value = 1171.6666
$ws.Range("K136").Value = 3514.9998
$ws.Range("M136").Value = -964.9998000000001
